# Final version for handing in report:
#  - fill in the R index column (A2:A24)
#  - correct / fill in NOTES (elastic/grain) column, and fill the missing
#    PREDICTION ("x") column for the rows that never got one
#  - add RADIANS() formulas for the rows that were missing THETA [RAD]
#  - add the "0.524…." annotation in D23
#  - move the active selection to K13 (and drop the old frozen/top-left scroll position)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- R index column (A2:A24) : 1..23 ---------------------------------------
for ($r = 2; $r -le 24; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# --- NOTES column (E) corrections / additions --------------------------------
$notes = @{
    2  = "elastic"
    10 = "grain"
    11 = "elastic"
    15 = "grain"
    16 = "elastic"
    17 = "elastic"
    18 = "grain"
    19 = "elastic"
    20 = "elastic"
    21 = "elastic"
    22 = "elastic"
    23 = "grain"
    24 = "elastic"
}
foreach ($r in $notes.Keys) {
    $ws.Cells.Item($r, 5).Value = $notes[$r]
}

# --- PREDICTION column (G) : rows 15-24 never had a value, fill with "x" ----
for ($r = 15; $r -le 24; $r++) {
    $ws.Cells.Item($r, 7).Value = "x"
}

# --- Newly-computed THETA [RAD] formulas (D17:D20) ---------------------------
$ws.Range("D17").Formula = "=RADIANS(C17)"
$ws.Range("D18").Formula = "=RADIANS(C18)"
$ws.Range("D19").Formula = "=RADIANS(C19)"
$ws.Range("D20").Formula = "=RADIANS(C20)"

# --- D23 gets a text annotation instead of a numeric THETA [RAD] ------------
$ws.Range("D23").Value = "0.524…."

# --- view: scroll/selection now sits on K13, no pinned topLeftCell ----------
$ws.Range("K13").Select()

Write-Host "edits applied"
